$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: updated timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 becomes Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the old duplicate "Contact" row (row 11), shifting rows 12-15 up
$ws.Rows.Item(11).Delete()
